$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before L to make room for the "Media" sector,
# shifting the former column L ("Fixed Income") to M.
$ws.Columns("L").Insert()

# New header for the inserted column
$ws.Range("L1").Value = "Media"

# Refresh all the sector-performance figures (Fall 2020 trade data refresh)
$row = New-Object 'object[,]' 1,12
$row[0,0] = -1.357
$row[0,1] = -0.418
$row[0,2] = -13.513
$row[0,3] = -2.211
$row[0,4] = -15.323
$row[0,5] = -0.671
$row[0,6] = -1.063
$row[0,7] = -2.419
$row[0,8] = -6.62
$row[0,9] = -3.21
$row[0,10] = -5.118
$row[0,11] = 0.857
$ws.Range("B2:M2").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = 2.161
$row[0,1] = 14.271
$row[0,2] = -16.967
$row[0,3] = 0.261
$row[0,4] = -8.663
$row[0,5] = 12.287
$row[0,6] = 5.043
$row[0,7] = 18.754
$row[0,8] = -0.296
$row[0,9] = 7.892
$row[0,10] = -4.058
$row[0,11] = 5.449
$ws.Range("B3:M3").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -12.295
$row[0,1] = 9.16
$row[0,2] = -34.716
$row[0,3] = -5.528
$row[0,4] = -14.427
$row[0,5] = 4.975
$row[0,6] = -32.25
$row[0,7] = -7.425
$row[0,8] = -9.298
$row[0,9] = 24.154
$row[0,10] = -6.717
$row[0,11] = -2.206
$ws.Range("B4:M4").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -8.537000000000001
$row[0,1] = 14.316
$row[0,2] = -33.963
$row[0,3] = 1.783
$row[0,4] = -14.32
$row[0,5] = 11.84
$row[0,6] = -28.917
$row[0,7] = -4.489
$row[0,8] = -1.983
$row[0,9] = 30.055
$row[0,10] = -0.755
$row[0,11] = 5.059
$ws.Range("B5:M5").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -8.446
$row[0,1] = 40.174
$row[0,2] = -40.337
$row[0,3] = 58.315
$row[0,4] = -25.232
$row[0,5] = 36.093
$row[0,6] = -29.178
$row[0,7] = 18.723
$row[0,8] = 5.33
$row[0,9] = 119.843
$row[0,10] = 54.679
$row[0,11] = -8.837
$ws.Range("B6:M6").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -1.767
$row[0,1] = 7.06
$row[0,2] = -9.906000000000001
$row[0,3] = 9.724
$row[0,4] = -5.704
$row[0,5] = 6.423
$row[0,6] = -6.731
$row[0,7] = 3.527
$row[0,8] = 1.054
$row[0,9] = 17.248
$row[0,10] = 9.210000000000001
$row[0,11] = -1.851
$ws.Range("B7:M7").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = 0.85
$row[0,1] = 0.978
$row[0,2] = 1.075
$row[0,3] = 0.971
$row[0,4] = 1.03
$row[0,5] = 1.02
$row[0,6] = 1.119
$row[0,7] = 0.881
$row[0,8] = 0.893
$row[0,9] = 1.063
$row[0,10] = 0.952
$row[0,11] = 0.598
$ws.Range("B8:M8").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -8.669
$row[0,1] = -0.781
$row[0,2] = -18.459
$row[0,3] = 1.933
$row[0,4] = -13.926
$row[0,5] = -1.722
$row[0,6] = -15.601
$row[0,7] = -3.601
$row[0,8] = -6.16
$row[0,9] = 8.788
$row[0,10] = 1.564
$row[0,11] = -6.908
$ws.Range("B9:M9").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -0.128
$row[0,1] = 0.267
$row[0,2] = -0.319
$row[0,3] = 0.369
$row[0,4] = -0.201
$row[0,5] = 0.237
$row[0,6] = -0.27
$row[0,7] = 0.113
$row[0,8] = 0.02
$row[0,9] = 0.577
$row[0,10] = 0.392
$row[0,11] = -0.204
$ws.Range("B10:M10").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -0.029
$row[0,1] = 0.065
$row[0,2] = -0.098
$row[0,3] = 0.093
$row[0,4] = -0.062
$row[0,5] = 0.056
$row[0,6] = -0.066
$row[0,7] = 0.032
$row[0,8] = 0.004
$row[0,9] = 0.156
$row[0,10] = 0.09
$row[0,11] = -0.042
$ws.Range("B11:M11").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = -46.362
$row[0,1] = -44.163
$row[0,2] = -66.36499999999999
$row[0,3] = -44.101
$row[0,4] = -53.669
$row[0,5] = -36.075
$row[0,6] = -50.385
$row[0,7] = -42.377
$row[0,8] = -43.677
$row[0,9] = -39.753
$row[0,10] = -35.67
$row[0,11] = -34.782
$ws.Range("B12:M12").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = 19.118
$row[0,1] = 23.976
$row[0,2] = 33.217
$row[0,3] = 24.559
$row[0,4] = 31.731
$row[0,5] = 24.251
$row[0,6] = 27.382
$row[0,7] = 25.248
$row[0,8] = 19.048
$row[0,9] = 28.76
$row[0,10] = 21.804
$row[0,11] = 12.393
$ws.Range("B13:M13").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = 0.45
$row[0,1] = 0.498
$row[0,2] = 0.471
$row[0,3] = 0.422
$row[0,4] = 0.393
$row[0,5] = 0.512
$row[0,6] = 0.471
$row[0,7] = 0.326
$row[0,8] = 0.607
$row[0,9] = 0.593
$row[0,10] = 0.674
$row[0,11] = 0.296
$ws.Range("B14:M14").Value = $row

$row = New-Object 'object[,]' 1,12
$row[0,0] = 6.902
$row[0,1] = 7.84
$row[0,2] = 8.552
$row[0,3] = 7.791
$row[0,4] = 8.221
$row[0,5] = 8.145
$row[0,6] = 8.869999999999999
$row[0,7] = 7.128
$row[0,8] = 7.214
$row[0,9] = 8.460000000000001
$row[0,10] = 7.646
$row[0,11] = 5.057
$ws.Range("B15:M15").Value = $row
